# Update "想去人数" (number of people interested) counts that changed
# between the previous gh-pages data export and the one generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6683
$wsExhibit.Range("F3").Value = 44
$wsExhibit.Range("F4").Value = 194
$wsExhibit.Range("F5").Value = 1044
$wsExhibit.Range("F6").Value = 137

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 7

# --- Sheet "全部类型" (All types, combined) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6683
$wsAll.Range("F3").Value = 44
$wsAll.Range("F4").Value = 194
$wsAll.Range("F5").Value = 1044
$wsAll.Range("F6").Value = 137
$wsAll.Range("F7").Value = 7
